# Scheduled-runner refresh of market-price-derived columns (H..N) across the
# per-job-class Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# All touched cells are plain numeric literals (no formulas in this workbook),
# so each changed value is written directly. A couple of rows also had a
# trailing N (or M/N) cell appear/disappear entirely between the two market
# snapshots, so those are added via a plain Value assignment (Excel creates
# the cell) or removed via ClearContents() (Excel drops the cell) instead of
# just being zeroed out.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 164.28572
$ws.Range("I12").Value = 225
$ws.Range("J12").Value = 83.333336
$ws.Range("K12").Value = 225
$ws.Range("L12").Value = 83.333336
$ws.Range("M12").Value = -55
$ws.Range("N12").Value = -423.333336

$ws.Range("H45").Value = 1511.3334

$ws.Range("H48").Value = 750
$ws.Range("J48").Value = 750
$ws.Range("L48").Value = 2250
$ws.Range("N48").Value = -2834

$ws.Range("H56").Value = 750
$ws.Range("J56").Value = 750
$ws.Range("L56").Value = 2250
$ws.Range("N56").Value = -3318

$ws.Range("H69").Value = 3735.4194
$ws.Range("I69").Value = 3710.818
$ws.Range("J69").Value = 3748.95
$ws.Range("K69").Value = 11132.454
$ws.Range("L69").Value = 11246.85
$ws.Range("M69").Value = -10258.454
$ws.Range("N69").Value = -12994.85

$ws.Range("H72").Value = 3735.4194
$ws.Range("I72").Value = 3710.818
$ws.Range("J72").Value = 3748.95
$ws.Range("K72").Value = 33397.362
$ws.Range("L72").Value = 33740.55
$ws.Range("M72").Value = -29029.362
$ws.Range("N72").Value = -42476.55

$ws.Range("H74").Value = 4484.8887
$ws.Range("I74").Value = 4138
$ws.Range("J74").Value = 4918.5
$ws.Range("K74").Value = 4138
$ws.Range("L74").Value = 4918.5
$ws.Range("M74").Value = -3202
$ws.Range("N74").Value = -6790.5

$ws.Range("H77").Value = 4484.8887
$ws.Range("I77").Value = 4138
$ws.Range("J77").Value = 4918.5
$ws.Range("K77").Value = 20690
$ws.Range("L77").Value = 24592.5
$ws.Range("M77").Value = -16010
$ws.Range("N77").Value = -33952.5

$ws.Range("H129").Value = 1055.9615
$ws.Range("J129").Value = 1206.975
$ws.Range("L129").Value = 3620.925
$ws.Range("N129").Value = -13620.925

$ws.Range("H138").Value = 2391.411
$ws.Range("I138").Value = 1442.1389
$ws.Range("J138").Value = 3315.027
$ws.Range("K138").Value = 4326.4167
$ws.Range("L138").Value = 9945.081
$ws.Range("M138").Value = 813.5833000000002
$ws.Range("N138").Value = -20225.081

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1618.6666
$ws.Range("I45").Value = 1820.6666
$ws.Range("J45").Value = 1416.6666
$ws.Range("K45").Value = 1820.6666
$ws.Range("L45").Value = 1416.6666
$ws.Range("M45").Value = -1443.6666
$ws.Range("N45").Value = -2170.6666

$ws.Range("H61").Value = 3091.7693
$ws.Range("I61").Value = 2526.2222
$ws.Range("J61").Value = 4364.25
$ws.Range("K61").Value = 2526.2222
$ws.Range("L61").Value = 4364.25
$ws.Range("M61").Value = -2314.2222
$ws.Range("N61").Value = -4788.25

$ws.Range("H76").Value = 16818.5
$ws.Range("J76").Value = 16818.5
$ws.Range("L76").Value = 16818.5
$ws.Range("N76").Value = -17494.5

$ws.Range("H79").Value = 16818.5
$ws.Range("J79").Value = 16818.5
$ws.Range("L79").Value = 16818.5
$ws.Range("N79").Value = -19158.5

$ws.Range("H132").Value = 3273.288
$ws.Range("I132").Value = 2229.4666
$ws.Range("K132").Value = 6688.399800000001
$ws.Range("M132").Value = -4158.399800000001

$ws.Range("H136").Value = 3091.7693
$ws.Range("I136").Value = 2526.2222
$ws.Range("J136").Value = 4364.25
$ws.Range("K136").Value = 7578.6666
$ws.Range("L136").Value = 13092.75
$ws.Range("M136").Value = -5028.6666
$ws.Range("N136").Value = -18192.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1190.6666
$ws.Range("I80").Value = 2740.7144
$ws.Range("J80").Value = 204.27272
$ws.Range("K80").Value = 2740.7144
$ws.Range("L80").Value = 204.27272
$ws.Range("M80").Value = -1742.7144
$ws.Range("N80").Value = -2200.27272

$ws.Range("H83").Value = 1190.6666
$ws.Range("I83").Value = 2740.7144
$ws.Range("J83").Value = 204.27272
$ws.Range("K83").Value = 13703.572
$ws.Range("L83").Value = 1021.3636
$ws.Range("M83").Value = -8711.572
$ws.Range("N83").Value = -11005.3636

$ws.Range("H86").Value = 3497.6667
$ws.Range("I86").Value = 3497.6667
$ws.Range("K86").Value = 3497.6667
$ws.Range("M86").Value = -2374.6667

$ws.Range("H89").Value = 3497.6667
$ws.Range("I89").Value = 3497.6667
$ws.Range("K89").Value = 17488.3335
$ws.Range("M89").Value = -11872.3335

$ws.Range("H105").Value = 2833.3333
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753

$ws.Range("H124").Value = 37723.332
$ws.Range("J124").Value = 37723.332
$ws.Range("L124").Value = 37723.332
$ws.Range("N124").Value = -47543.332

$ws.Range("H132").Value = 42866.43
$ws.Range("J132").Value = 42866.43
$ws.Range("L132").Value = 42866.43
$ws.Range("N132").Value = -52986.43

$ws.Range("H134").Value = 2531.96
$ws.Range("I134").Value = 2252.3333
$ws.Range("K134").Value = 6756.999899999999
$ws.Range("M134").Value = -4221.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4387834.5
$ws.Range("I132").Value = 1333
$ws.Range("J132").Value = 12823414
$ws.Range("K132").Value = 3999
$ws.Range("L132").Value = 38470242
$ws.Range("M132").Value = -1469
$ws.Range("N132").Value = -38475302

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 37.4375
$ws.Range("I2").Value = 19.428572
$ws.Range("J2").Value = 51.444443
$ws.Range("K2").Value = 116.571432
$ws.Range("L2").Value = 308.666658
$ws.Range("M2").Value = -3.571431999999987
$ws.Range("N2").Value = -534.666658

$ws.Range("H76").Value = 3989
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3989
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H115").Value = 5453.143
$ws.Range("I115").Value = 5003.3335
$ws.Range("J115").Value = 5790.5
$ws.Range("K115").Value = 15010.0005
$ws.Range("L115").Value = 17371.5
$ws.Range("M115").Value = -13835.0005
$ws.Range("N115").Value = -19721.5

$ws.Range("H119").Value = 2965.8
$ws.Range("I119").Value = 2965.8
$ws.Range("K119").Value = 8897.400000000001
$ws.Range("M119").Value = -4059.400000000001

$ws.Range("H140").Value = 2101.6667
$ws.Range("J140").Value = 2211.4285
$ws.Range("L140").Value = 6634.2855
$ws.Range("N140").Value = -16994.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1757.2222
$ws.Range("I102").Value = 1676.3334
$ws.Range("J102").Value = 1919
$ws.Range("K102").Value = 1676.3334
$ws.Range("L102").Value = 1919
$ws.Range("M102").Value = -54.33339999999998
$ws.Range("N102").Value = -5163

$ws.Range("H116").Value = 39999
$ws.Range("J116").Value = 39999
$ws.Range("L116").Value = 39999
$ws.Range("N116").Value = -49177

$ws.Range("H122").Value = 2204.7144
$ws.Range("I122").Value = 2108.25
$ws.Range("K122").Value = 6324.75
$ws.Range("M122").Value = -3874.75

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3326.4285
$ws.Range("I132").Value = 2867.2222
$ws.Range("K132").Value = 8601.6666
$ws.Range("M132").Value = -6071.6666

$ws.Range("H133").Value = 47731.668
$ws.Range("J133").Value = 47731.668
$ws.Range("L133").Value = 47731.668
$ws.Range("N133").Value = -57851.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4167929.2
$ws.Range("I136").Value = 1109.037
$ws.Range("J136").Value = 12822094
$ws.Range("K136").Value = 3327.111
$ws.Range("L136").Value = 38466282
$ws.Range("M136").Value = -777.1109999999999
$ws.Range("N136").Value = -38471382

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2925557.5
$ws.Range("I132").Value = 1484.8049
$ws.Range("J132").Value = 10418494
$ws.Range("K132").Value = 4454.4147
$ws.Range("L132").Value = 31255482
$ws.Range("M132").Value = -1924.4147
$ws.Range("N132").Value = -31260542
